$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells that should no longer hold a value
$ws.Range("D4").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("D8").ClearContents()

# Update the active selection to D4
$ws.Range("D4").Select()
